$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.6542786666666667
$ws.Cells.Item(2, 8).Value = 1.962836
$ws.Cells.Item(2, 9).Value = 0.3193006097963691
$ws.Cells.Item(2, 10).Value = 0.3193006097963691
$ws.Cells.Item(2, 13).Value = 209.26237
$ws.Cells.Item(2, 14).Value = 627.78711
$ws.Cells.Item(2, 15).Value = 0.8127157202241573
$ws.Cells.Item(2, 16).Value = 0.8127157202241573
$ws.Cells.Item(2, 17).Value = 136.9159044271067
$ws.Cells.Item(2, 18).Value = 1232.24313984396
$ws.Cells.Item(2, 19).Value = 0.2595006250586687
$ws.Cells.Item(2, 20).Value = 0.2595006250586687

# Row 3
$ws.Cells.Item(3, 7).Value = 0.6542786666666667
$ws.Cells.Item(3, 8).Value = 1.962836
$ws.Cells.Item(3, 9).Value = 0.3193006097963691
$ws.Cells.Item(3, 10).Value = 0.3193006097963691
$ws.Cells.Item(3, 13).Value = 0.9848756666666668
$ws.Cells.Item(3, 14).Value = 2.954627
$ws.Cells.Item(3, 15).Value = 0.003824977881910862
$ws.Cells.Item(3, 16).Value = 0.003824977881910862
$ws.Cells.Item(3, 17).Value = 0.6443831380191112
$ws.Cells.Item(3, 18).Value = 5.799448242172001
$ws.Cells.Item(3, 19).Value = 0.001221317770151763
$ws.Cells.Item(3, 20).Value = 0.001221317770151763

# Row 4
$ws.Cells.Item(4, 7).Value = 0.6542786666666667
$ws.Cells.Item(4, 8).Value = 1.962836
$ws.Cells.Item(4, 9).Value = 0.3193006097963691
$ws.Cells.Item(4, 10).Value = 0.3193006097963691
$ws.Cells.Item(4, 13).Value = 1.763846666666667
$ws.Cells.Item(4, 14).Value = 5.291539999999999
$ws.Cells.Item(4, 15).Value = 0.006850280411451801
$ws.Cells.Item(4, 16).Value = 0.006850280411451801
$ws.Cells.Item(4, 17).Value = 1.154047245271111
$ws.Cells.Item(4, 18).Value = 10.38642520744
$ws.Cells.Item(4, 19).Value = 0.002187298712652682
$ws.Cells.Item(4, 20).Value = 0.002187298712652682

# Row 5
$ws.Cells.Item(5, 7).Value = 0.6542786666666667
$ws.Cells.Item(5, 8).Value = 1.962836
$ws.Cells.Item(5, 9).Value = 0.3193006097963691
$ws.Cells.Item(5, 10).Value = 0.3193006097963691
$ws.Cells.Item(5, 13).Value = 45.474231
$ws.Cells.Item(5, 14).Value = 136.422693
$ws.Cells.Item(5, 15).Value = 0.1766090214824801
$ws.Cells.Item(5, 16).Value = 0.1766090214824801
$ws.Cells.Item(5, 17).Value = 29.752819226372
$ws.Cells.Item(5, 18).Value = 267.7753730373479
$ws.Cells.Item(5, 19).Value = 0.05639136825489594
$ws.Cells.Item(5, 20).Value = 0.05639136825489594

# Row 6
$ws.Cells.Item(6, 9).Value = 0.4124821994964292
$ws.Cells.Item(6, 10).Value = 0.4124821994964292
$ws.Cells.Item(6, 13).Value = 209.26237
$ws.Cells.Item(6, 14).Value = 627.78711
$ws.Cells.Item(6, 15).Value = 0.8127157202241573
$ws.Cells.Item(6, 16).Value = 0.8127157202241573
$ws.Cells.Item(6, 17).Value = 176.87211258429
$ws.Cells.Item(6, 18).Value = 1591.84901325861
$ws.Cells.Item(6, 19).Value = 0.335230767843385
$ws.Cells.Item(6, 20).Value = 0.335230767843385

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4124821994964292
$ws.Cells.Item(7, 10).Value = 0.4124821994964292
$ws.Cells.Item(7, 13).Value = 0.9848756666666668
$ws.Cells.Item(7, 14).Value = 2.954627
$ws.Cells.Item(7, 15).Value = 0.003824977881910862
$ws.Cells.Item(7, 16).Value = 0.003824977881910862
$ws.Cells.Item(7, 17).Value = 0.8324336563530002
$ws.Cells.Item(7, 18).Value = 7.491902907177002
$ws.Cells.Item(7, 19).Value = 0.001577735289755785
$ws.Cells.Item(7, 20).Value = 0.001577735289755785

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4124821994964292
$ws.Cells.Item(8, 10).Value = 0.4124821994964292
$ws.Cells.Item(8, 13).Value = 1.763846666666667
$ws.Cells.Item(8, 14).Value = 5.291539999999999
$ws.Cells.Item(8, 15).Value = 0.006850280411451801
$ws.Cells.Item(8, 16).Value = 0.006850280411451801
$ws.Cells.Item(8, 17).Value = 1.49083318806
$ws.Cells.Item(8, 18).Value = 13.41749869254
$ws.Cells.Item(8, 19).Value = 0.002825618731282942
$ws.Cells.Item(8, 20).Value = 0.002825618731282942

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4124821994964292
$ws.Cells.Item(9, 10).Value = 0.4124821994964292
$ws.Cells.Item(9, 13).Value = 45.474231
$ws.Cells.Item(9, 14).Value = 136.422693
$ws.Cells.Item(9, 15).Value = 0.1766090214824801
$ws.Cells.Item(9, 16).Value = 0.1766090214824801
$ws.Cells.Item(9, 17).Value = 38.435593103127
$ws.Cells.Item(9, 18).Value = 345.9203379281429
$ws.Cells.Item(9, 19).Value = 0.0728480776320055
$ws.Cells.Item(9, 20).Value = 0.0728480776320055

# Row 10
$ws.Cells.Item(10, 7).Value = 0.5380133333333333
$ws.Cells.Item(10, 8).Value = 1.61404
$ws.Cells.Item(10, 9).Value = 0.2625608844731457
$ws.Cells.Item(10, 10).Value = 0.2625608844731457
$ws.Cells.Item(10, 13).Value = 209.26237
$ws.Cells.Item(10, 14).Value = 627.78711
$ws.Cells.Item(10, 15).Value = 0.8127157202241573
$ws.Cells.Item(10, 16).Value = 0.8127157202241573
$ws.Cells.Item(10, 17).Value = 112.5859452249333
$ws.Cells.Item(10, 18).Value = 1013.2735070244
$ws.Cells.Item(10, 19).Value = 0.2133873583272844
$ws.Cells.Item(10, 20).Value = 0.2133873583272844

# Row 11
$ws.Cells.Item(11, 7).Value = 0.5380133333333333
$ws.Cells.Item(11, 8).Value = 1.61404
$ws.Cells.Item(11, 9).Value = 0.2625608844731457
$ws.Cells.Item(11, 10).Value = 0.2625608844731457
$ws.Cells.Item(11, 13).Value = 0.9848756666666668
$ws.Cells.Item(11, 14).Value = 2.954627
$ws.Cells.Item(11, 15).Value = 0.003824977881910862
$ws.Cells.Item(11, 16).Value = 0.003824977881910862
$ws.Cells.Item(11, 17).Value = 0.5298762403422224
$ws.Cells.Item(11, 18).Value = 4.76888616308
$ws.Cells.Item(11, 19).Value = 0.001004289575764736
$ws.Cells.Item(11, 20).Value = 0.001004289575764736

# Row 12
$ws.Cells.Item(12, 7).Value = 0.5380133333333333
$ws.Cells.Item(12, 8).Value = 1.61404
$ws.Cells.Item(12, 9).Value = 0.2625608844731457
$ws.Cells.Item(12, 10).Value = 0.2625608844731457
$ws.Cells.Item(12, 13).Value = 1.763846666666667
$ws.Cells.Item(12, 14).Value = 5.291539999999999
$ws.Cells.Item(12, 15).Value = 0.006850280411451801
$ws.Cells.Item(12, 16).Value = 0.006850280411451801
$ws.Cells.Item(12, 17).Value = 0.9489730246222222
$ws.Cells.Item(12, 18).Value = 8.540757221599998
$ws.Cells.Item(12, 19).Value = 0.00179861568371985
$ws.Cells.Item(12, 20).Value = 0.00179861568371985

# Row 13
$ws.Cells.Item(13, 7).Value = 0.5380133333333333
$ws.Cells.Item(13, 8).Value = 1.61404
$ws.Cells.Item(13, 9).Value = 0.2625608844731457
$ws.Cells.Item(13, 10).Value = 0.2625608844731457
$ws.Cells.Item(13, 13).Value = 45.474231
$ws.Cells.Item(13, 14).Value = 136.422693
$ws.Cells.Item(13, 15).Value = 0.1766090214824801
$ws.Cells.Item(13, 16).Value = 0.1766090214824801
$ws.Cells.Item(13, 17).Value = 24.46574260108
$ws.Cells.Item(13, 18).Value = 220.19168340972
$ws.Cells.Item(13, 19).Value = 0.04637062088637677
$ws.Cells.Item(13, 20).Value = 0.04637062088637677

# Row 14
$ws.Cells.Item(14, 7).Value = 0.01159033333333333
$ws.Cells.Item(14, 8).Value = 0.034771
$ws.Cells.Item(14, 9).Value = 0.005656306234056004
$ws.Cells.Item(14, 10).Value = 0.005656306234056004
$ws.Cells.Item(14, 13).Value = 209.26237
$ws.Cells.Item(14, 14).Value = 627.78711
$ws.Cells.Item(14, 15).Value = 0.8127157202241573
$ws.Cells.Item(14, 16).Value = 0.8127157202241573
$ws.Cells.Item(14, 17).Value = 2.425420622423334
$ws.Cells.Item(14, 18).Value = 21.82878560181
$ws.Cells.Item(14, 19).Value = 0.004596968994819216
$ws.Cells.Item(14, 20).Value = 0.004596968994819216

# Row 15
$ws.Cells.Item(15, 7).Value = 0.01159033333333333
$ws.Cells.Item(15, 8).Value = 0.034771
$ws.Cells.Item(15, 9).Value = 0.005656306234056004
$ws.Cells.Item(15, 10).Value = 0.005656306234056004
$ws.Cells.Item(15, 13).Value = 0.9848756666666668
$ws.Cells.Item(15, 14).Value = 2.954627
$ws.Cells.Item(15, 15).Value = 0.003824977881910862
$ws.Cells.Item(15, 16).Value = 0.003824977881910862
$ws.Cells.Item(15, 17).Value = 0.01141503726855556
$ws.Cells.Item(15, 18).Value = 0.102735335417
$ws.Cells.Item(15, 19).Value = 0.00002163524623857874
$ws.Cells.Item(15, 20).Value = 0.00002163524623857874

# Row 16
$ws.Cells.Item(16, 7).Value = 0.01159033333333333
$ws.Cells.Item(16, 8).Value = 0.034771
$ws.Cells.Item(16, 9).Value = 0.005656306234056004
$ws.Cells.Item(16, 10).Value = 0.005656306234056004
$ws.Cells.Item(16, 13).Value = 1.763846666666667
$ws.Cells.Item(16, 14).Value = 5.291539999999999
$ws.Cells.Item(16, 15).Value = 0.006850280411451801
$ws.Cells.Item(16, 16).Value = 0.006850280411451801
$ws.Cells.Item(16, 17).Value = 0.02044357081555555
$ws.Cells.Item(16, 18).Value = 0.18399213734
$ws.Cells.Item(16, 19).Value = 0.00003874728379632655
$ws.Cells.Item(16, 20).Value = 0.00003874728379632655

# Row 17
$ws.Cells.Item(17, 7).Value = 0.01159033333333333
$ws.Cells.Item(17, 8).Value = 0.034771
$ws.Cells.Item(17, 9).Value = 0.005656306234056004
$ws.Cells.Item(17, 10).Value = 0.005656306234056004
$ws.Cells.Item(17, 13).Value = 45.474231
$ws.Cells.Item(17, 14).Value = 136.422693
$ws.Cells.Item(17, 15).Value = 0.1766090214824801
$ws.Cells.Item(17, 16).Value = 0.1766090214824801
$ws.Cells.Item(17, 17).Value = 0.527061495367
$ws.Cells.Item(17, 18).Value = 4.743553458302999
$ws.Cells.Item(17, 19).Value = 0.0009989547092018828
$ws.Cells.Item(17, 20).Value = 0.0009989547092018828
